$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F ("想去人数") on sheets "展览" and "全部类型"
$updates = @{
    2  = 199
    3  = 241
    5  = 796
    7  = 5999
    8  = 41
    9  = 69
    11 = 52
    14 = 182
    15 = 383
    16 = 32
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
